$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("A9").Value = 112171810
$ws.Range("B9").Value = 90480
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 4769
$ws.Range("F9").Value = "Svavelriska"
$ws.Range("G9").Value = "Lactarius scrobiculatus"
$ws.Range("H9").Value = "(Scop.:Fr.) Fr."
$ws.Range("Q9").Value = 756486
$ws.Range("R9").Value = 7212020
# Row 10
$ws.Range("A10").Value = 112171787
$ws.Range("B10").Value = 86371
$ws.Range("E10").Value = 4412
$ws.Range("F10").Value = "Äggvaxskivling"
$ws.Range("G10").Value = "Hygrophorus karstenii"
$ws.Range("H10").Value = "Sacc. & Cub."
$ws.Range("Q10").Value = 756408
$ws.Range("R10").Value = 7211956
# Row 11
$ws.Range("A11").Value = 112171806
$ws.Range("B11").Value = 86371
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 4412
$ws.Range("F11").Value = "Äggvaxskivling"
$ws.Range("G11").Value = "Hygrophorus karstenii"
$ws.Range("H11").Value = "Sacc. & Cub."
$ws.Range("Q11").Value = 756477
$ws.Range("R11").Value = 7212031
# Row 12
$ws.Range("A12").Value = 112171785
$ws.Range("B12").Value = 78677
$ws.Range("D12").Value = "LC"
$ws.Range("E12").Value = 229748
$ws.Range("F12").Value = "Gytterlav"
$ws.Range("G12").Value = "Protopannaria pezizoides"
$ws.Range("H12").Value = "(Weber) P.M.Jørg. & S.Ekman"
$ws.Range("Q12").Value = 756412
$ws.Range("R12").Value = 7211954
$ws.Range("AC12").Value = "på berg"
# Row 13
$ws.Range("A13").Value = 112171795
$ws.Range("B13").Value = 77650
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = "Garnlav"
$ws.Range("G13").Value = "Alectoria sarmentosa"
$ws.Range("Q13").Value = 756378
$ws.Range("R13").Value = 7212050
# Row 14
$ws.Range("A14").Value = 112171812
$ws.Range("B14").Value = 78746
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 6463
$ws.Range("F14").Value = "Bårdlav"
$ws.Range("G14").Value = "Nephroma parile"
$ws.Range("H14").Value = "(Ach.) Ach."
$ws.Range("Q14").Value = 756485
$ws.Range("R14").Value = 7212023
$ws.Range("AJ14").Value = "sälg"
$ws.Range("AK14").Value = "Salix caprea"
$ws.Range("AO14").Value = "Salix caprea"
# Row 15
$ws.Range("A15").Value = 112171813
$ws.Range("B15").Value = 78713
$ws.Range("E15").Value = 6458
$ws.Range("F15").Value = "Lunglav"
$ws.Range("G15").Value = "Lobaria pulmonaria"
$ws.Range("H15").Value = "(L.) Hoffm."
$ws.Range("Q15").Value = 756485
$ws.Range("R15").Value = 7212023
$ws.Range("AJ15").Value = "sälg"
$ws.Range("AK15").Value = "Salix caprea"
$ws.Range("AO15").Value = "Salix caprea"
# Row 16
$ws.Range("A16").Value = 112171801
$ws.Range("B16").Value = 78713
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 6458
$ws.Range("F16").Value = "Lunglav"
$ws.Range("G16").Value = "Lobaria pulmonaria"
$ws.Range("H16").Value = "(L.) Hoffm."
$ws.Range("Q16").Value = 756448
$ws.Range("R16").Value = 7212052
# Row 17
$ws.Range("A17").Value = 112171798
$ws.Range("B17").Value = 78713
$ws.Range("D17").Value = "NT"
$ws.Range("E17").Value = 6458
$ws.Range("F17").Value = "Lunglav"
$ws.Range("G17").Value = "Lobaria pulmonaria"
$ws.Range("H17").Value = "(L.) Hoffm."
$ws.Range("Q17").Value = 756371
$ws.Range("R17").Value = 7212116
# Row 18
$ws.Range("A18").Value = 112171814
$ws.Range("B18").Value = 89499
$ws.Range("E18").Value = 112
$ws.Range("F18").Value = "Stjärntagging"
$ws.Range("G18").Value = "Asterodon ferruginosus"
$ws.Range("H18").Value = "Pat."
$ws.Range("Q18").Value = 756486
$ws.Range("R18").Value = 7212041
$ws.Range("AJ18").Value = "gran"
$ws.Range("AK18").Value = "Picea abies"
$ws.Range("AO18").Value = "Picea abies"
# Row 19
$ws.Range("A19").Value = 112171779
$ws.Range("B19").Value = 102192
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 222412
$ws.Range("F19").Value = "Tibast"
$ws.Range("G19").Value = "Daphne mezereum"
$ws.Range("H19").Value = "L."
$ws.Range("Q19").Value = 756291
$ws.Range("R19").Value = 7211892
# Row 20
$ws.Range("A20").Value = 112171776
$ws.Range("B20").Value = 85448
$ws.Range("E20").Value = 3739
$ws.Range("F20").Value = "Persiljespindling"
$ws.Range("G20").Value = "Cortinarius sulfurinus"
$ws.Range("H20").Value = "Quél."
$ws.Range("Q20").Value = 756261
$ws.Range("R20").Value = 7211953
# Row 21
$ws.Range("A21").Value = 112171792
$ws.Range("B21").Value = 85401
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 249228
$ws.Range("F21").Value = "Barrfagerspindling"
$ws.Range("G21").Value = "Cortinarius piceae"
$ws.Range("H21").Value = "Frøslev, T.S.Jeppesen & Brandrud"
$ws.Range("Q21").Value = 756395
$ws.Range("R21").Value = 7211974
# Row 22
$ws.Range("A22").Value = 112171788
$ws.Range("B22").Value = 78713
$ws.Range("D22").Value = "NT"
$ws.Range("E22").Value = 6458
$ws.Range("F22").Value = "Lunglav"
$ws.Range("G22").Value = "Lobaria pulmonaria"
$ws.Range("H22").Value = "(L.) Hoffm."
$ws.Range("Q22").Value = 756401
$ws.Range("R22").Value = 7211954
$ws.Range("AJ22").Value = "sälg"
$ws.Range("AK22").Value = "Salix caprea"
$ws.Range("AO22").Value = "Salix caprea"

# Clear cells no longer populated
$ws.Range("AJ12").ClearContents()
$ws.Range("AK12").ClearContents()
$ws.Range("AO12").ClearContents()
$ws.Range("AJ13").ClearContents()
$ws.Range("AK13").ClearContents()
$ws.Range("AO13").ClearContents()
$ws.Range("AC16").ClearContents()
$ws.Range("AJ19").ClearContents()
$ws.Range("AK19").ClearContents()
$ws.Range("AO19").ClearContents()
$ws.Range("AJ20").ClearContents()
$ws.Range("AK20").ClearContents()
$ws.Range("AO20").ClearContents()
